$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Insert two new rows above the existing row 31 (old "DRAIAM102" row),
# shifting the current rows 31-36 down to 33-38.
$ws.Rows.Item(31).Resize(2).Insert()

# Fill column-by-column (A31, A32, B31, B32, C31, C32, D31, D32) to match
# the authoring order of the new shared-string table entries.
$ws.Cells.Item(31, 1).Value = "DRAIAM072"
$ws.Cells.Item(32, 1).Value = "DRAIAM073"

$ws.Cells.Item(31, 2).Value = "OPQA-5155"
$ws.Cells.Item(32, 2).Value = "OPQA-5228"

$ws.Cells.Item(31, 3).Value = "Verify that the user is able to manually select any particular country and see contact details associated with it as per IPA Customer Care Contact Details.doc and DRA Customer Care Contact Details.doc documents"
$ws.Cells.Item(32, 3).Value = "Verify that the page shall be accessible in both an authenticated and a non-authenticated state"

$ws.Cells.Item(31, 4).Value = "Y"
$ws.Cells.Item(32, 4).Value = "Y"

# Match the style/row-height used by neighboring rows.
$ws.Rows.Item(31).RowHeight = 30
$ws.Rows.Item(32).RowHeight = 15

$ws.Range("A31:E31").Style = $ws.Range("A33:E33").Style
$ws.Cells.Item(31, 4).Style = $ws.Cells.Item(30, 4).Style

$ws.Range("A32:E32").Style = $ws.Range("A33:E33").Style
$ws.Cells.Item(32, 4).Style = $ws.Cells.Item(30, 4).Style

# Update the view/selection to mirror the edited state.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("K33").Select()
